# feat: add 2022-Q4 data
#
# - "总计" sheet: a new "2022-Q4" row is inserted at the top of the data
#   (count=7, value=0.27); the existing "2022-Q2" / "2022-Q1" rows shift
#   down one position each (their own data is unchanged).
# - A brand-new "2022-Q4" worksheet is added (positioned right after the
#   "总计" sheet, before "2022-Q2") holding the per-fund breakdown for the
#   new quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" (summary) sheet: add the 2022-Q4 row, keep the
#    older quarters but shift them down.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Clone the style of row 3 (s="2": centered / bold / bordered) onto the
# brand-new row 4 so the new "2022-Q1" entry keeps the same look as the
# other data rows.
$summary.Range("A3").Copy()
$summary.Range("A4").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 7
$summary.Range("D2").Value = 0.27

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 3
$summary.Range("D3").Value = 0.2

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 3
$summary.Range("D4").Value = 0.32

# ---------------------------------------------------------------------
# 2) Add the new "2022-Q4" worksheet. Duplicate the "2022-Q2" sheet so
#    it starts out with identical column layout / header styling, drop
#    it in right before "2022-Q2", rename it, then replace the data.
# ---------------------------------------------------------------------
$sheets = $wb.Worksheets
$quarterTemplate = $sheets.Item("2022-Q2")
$quarterTemplate.Copy($quarterTemplate)
$new = $sheets.Item("2022-Q2 (2)")
$new.Name = "2022-Q4"

# The template only has 4 rows (1 header + 3 funds); we need 8 (1 header
# + 7 funds). Extend the formatting (text columns, A-column style, etc.)
# of row 2 down through rows 5-8 before writing the extra fund rows.
$new.Range("A2:H2").Copy()
$new.Range("A5:H8").PasteSpecial(-4122)

# Columns B, D, E, F, G hold numeric-looking values that must stay text
# (fund codes with leading zeros, percentages formatted as "xx.xx", …),
# same convention as the other quarter sheets in this workbook.
$textCols = @("B", "D", "E", "F", "G")
foreach ($col in $textCols) {
    $new.Range($col + "2:" + $col + "8").NumberFormat = "@"
}

$new.Range("A2").Value = 0
$new.Range("B2").Value = "006234"
$new.Range("C2").Value = "万家汽车新趋势混合C"
$new.Range("D2").Value = "4.06"
$new.Range("E2").Value = "90.27"
$new.Range("F2").Value = "3.37"
$new.Range("G2").Value = "0.1368"
$new.Range("H2").Value = 9

$new.Range("A3").Value = 1
$new.Range("B3").Value = "014339"
$new.Range("C3").Value = "长江智能制造混合A"
$new.Range("D3").Value = "2.51"
$new.Range("E3").Value = "83.09"
$new.Range("F3").Value = "2.84"
$new.Range("G3").Value = "0.0713"
$new.Range("H3").Value = 8

$new.Range("A4").Value = 2
$new.Range("B4").Value = "006233"
$new.Range("C4").Value = "万家汽车新趋势混合A"
$new.Range("D4").Value = "1.65"
$new.Range("E4").Value = "90.27"
$new.Range("F4").Value = "3.37"
$new.Range("G4").Value = "0.0556"
$new.Range("H4").Value = 9

$new.Range("A5").Value = 3
$new.Range("B5").Value = "001318"
$new.Range("C5").Value = "东方新策略灵活配置混合A"
$new.Range("D5").Value = "0.39"
$new.Range("E5").Value = "36.37"
$new.Range("F5").Value = "1.32"
$new.Range("G5").Value = "0.0051"
$new.Range("H5").Value = 2

$new.Range("A6").Value = 4
$new.Range("B6").Value = "014340"
$new.Range("C6").Value = "长江智能制造混合C"
$new.Range("D6").Value = "0.09"
$new.Range("E6").Value = "83.09"
$new.Range("F6").Value = "2.84"
$new.Range("G6").Value = "0.0026"
$new.Range("H6").Value = 8

$new.Range("A7").Value = 5
$new.Range("B7").Value = "400020"
$new.Range("C7").Value = "东方成长回报平衡混合"
$new.Range("D7").Value = "0.15"
$new.Range("E7").Value = "45.42"
$new.Range("F7").Value = "1.48"
$new.Range("G7").Value = "0.0022"
$new.Range("H7").Value = 4

$new.Range("A8").Value = 6
$new.Range("B8").Value = "002060"
$new.Range("C8").Value = "东方新策略灵活配置混合C"
$new.Range("D8").Value = "0.03"
$new.Range("E8").Value = "36.37"
$new.Range("F8").Value = "1.32"
$new.Range("G8").Value = "0.0004"
$new.Range("H8").Value = 2
